$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 45981
$ws.Range("B2").Value = 10485.9637300029
$ws.Range("C2").Value = 9805.38651934963
$ws.Range("D2").Value = 11259.4
$ws.Range("E2").Value = 5324.11903850399
$ws.Range("F2").Value = 161.254398243901

# Row 3
$ws.Range("A3").Value = 45982
$ws.Range("B3").Value = 10455.1114409099
$ws.Range("C3").Value = 9419.00870283085
$ws.Range("D3").Value = 9003.4
$ws.Range("E3").Value = 6070.58961563777
$ws.Range("F3").Value = 270.258263269526

# Row 4
$ws.Range("A4").Value = 45983
$ws.Range("B4").Value = 4417.55882369736
$ws.Range("C4").Value = 7214.56389466841
$ws.Range("D4").Value = 9003.4
$ws.Range("E4").Value = 5989.80618230243
$ws.Range("F4").Value = 175.040419873785

# Row 5
$ws.Range("A5").Value = 45984
$ws.Range("B5").Value = 4551.30774402517
$ws.Range("C5").Value = 7551.70250462126
$ws.Range("D5").Value = 9003.4
$ws.Range("E5").Value = 6318.960247411
$ws.Range("F5").Value = 202.802614668011

# Row 6
$ws.Range("A6").Value = 45985
$ws.Range("B6").Value = 12485.6583594937
$ws.Range("C6").Value = 11067.6904838179
$ws.Range("D6").Value = 9003.4
$ws.Range("E6").Value = 7465.55050960317
$ws.Range("F6").Value = 397.076708059212

# Row 7
$ws.Range("A7").Value = 45986
$ws.Range("B7").Value = 10980.8696790314
$ws.Range("C7").Value = 10287.1212940648
$ws.Range("D7").Value = 9003.4
$ws.Range("E7").Value = 6409.24037742116
$ws.Range("F7").Value = 320.540069645248

# Row 8
$ws.Range("A8").Value = 45987
$ws.Range("B8").Value = 10980.8696790314
$ws.Range("C8").Value = 9979.53477819552
$ws.Range("D8").Value = 9003.4
$ws.Range("E8").Value = 6409.24037742116
$ws.Range("F8").Value = 307.723964817362

# Row 9
$ws.Range("A9").Value = 45988
$ws.Range("B9").Value = 10980.8696790314
$ws.Range("C9").Value = 10323.9141701013
$ws.Range("D9").Value = 9003.4
$ws.Range("E9").Value = 6409.24037742116
$ws.Range("F9").Value = 322.07310614677

# Row 10
$ws.Range("A10").Value = 45989
$ws.Range("B10").Value = 10980.8696790314
$ws.Range("C10").Value = 9893.14514510403
$ws.Range("D10").Value = 9003.4
$ws.Range("E10").Value = 6409.24037742116
$ws.Range("F10").Value = 304.124396771883

# Row 11
$ws.Range("A11").Value = 45990
$ws.Range("B11").Value = 4442.42485201967
$ws.Range("C11").Value = 6715.48138653909
$ws.Range("D11").Value = 9003.4
$ws.Range("E11").Value = 6033.02000427315
$ws.Range("F11").Value = 156.045891283843

# Row 12
$ws.Range("A12").Value = 45991
$ws.Range("B12").Value = 4286.74551609542
$ws.Range("C12").Value = 6551.76641709588
$ws.Range("D12").Value = 9003.4
$ws.Range("E12").Value = 6024.53142231133
$ws.Range("F12").Value = 148.870743308634

# Row 13
$ws.Range("A13").Value = 45992
$ws.Range("B13").Value = 11778.1169212215
$ws.Range("C13").Value = 10384.9443686614
$ws.Range("D13").Value = 8664.26
$ws.Range("E13").Value = 7489.73291767109
$ws.Range("F13").Value = 383.767386930521

# Row 14
$ws.Range("A14").Value = 45993
$ws.Range("B14").Value = 11778.1169212215
$ws.Range("C14").Value = 10535.5243695882
$ws.Range("D14").Value = 8664.26
$ws.Range("E14").Value = 7489.73291767109
$ws.Range("F14").Value = 390.041553635802

# Row 15
$ws.Range("A15").Value = 45994
$ws.Range("B15").Value = 11778.1169212215
$ws.Range("C15").Value = 10257.0116806562
$ws.Range("D15").Value = 8664.26
$ws.Range("E15").Value = 7489.73291767109
$ws.Range("F15").Value = 378.436858263637
